# Adds feedback notes (Darren / Matthew / Eve) plus an "Afterwards" wrap-up
# paragraph to the end of the "That Day" notes document, just before the
# trailing bookmark paragraph.

$d = $word.ActiveDocument

# The very last paragraph in the document is an (otherwise empty) paragraph
# that merely carries the _GoBack bookmark. We insert all of the new
# feedback paragraphs immediately before it, so the bookmark paragraph
# keeps trailing the document exactly as before.
$bookmarkPara = $d.Paragraphs.Last
$insertRange = $bookmarkPara.Range

$paragraphTexts = @(
    "Darren’s feedback:",
    "The villain is a weak character. The reader knows nothing about him, there is no relating to him. Unless the point was that Claire can’t keep running from her problems.",
    "The villain was really just a plot device to teach Claire a lesson. I will admit that the character that does occupy that space is one-dimensional and doesn’t really do very much. I could always flesh out the dynamic between him and Claire in another draft. I think he was just an oversight as I was so focused on Claire’s perspective and her interaction with her sister that the villain just got lost in the shuffle. He almost represents the hand of fate, coming into Claire’s life and forcing her to accept the truth.",
    "",
    "Matthew’s feedback:",
    "`tSibling is confusing.",
    "Plot is weak. Claire feels guiltier than she should about her involvement in the attack. Antagonist needs to have more of a personal connection to her if she is to feel as guilty as she does.",
    "`tClaire feels American.",
    "`t",
    "The plot being weak is a legitimate critique. ",
    "`t",
    "Eve’s feedback:",
    "`tThe use of the term ‘sibling’ is out of place because it is too formal.",
    "It gets difficult to discern whether the story is happening in the past or the present, particularly near to the end. ",
    "The term sibling is used to emphasize the emotional distance between Claire and Nicole, from Claire’s perspective at least. Will have to reread the essay to see how this happens and change it accordingly. It may be that Claire should, in the first half of the story, refer to Nicole as ‘sibling’ or ‘sister’ and then only call her by name later on. This will place emphasis on the emotional distance by creating a sense of abstraction. The reader is only allowed to see Nicole as a person when Claire does.",
    ""
)

# Paragraphs (1-based, within $paragraphTexts) that get a 0.5" (720 twips)
# left indent, matching <w:ind w:left="720"/> in the target markup.
$indentedIndexes = @(2, 7, 14)

# A trailing paragraph mark is required so that every entry in
# $paragraphTexts - including the last (empty) one - becomes its own
# paragraph, distinct from the pre-existing bookmark paragraph that the
# text is being inserted in front of.
$blockText = ([string]::Join("`r", $paragraphTexts)) + "`r"

# Remember where the newly-inserted paragraphs will start so we can revisit
# them afterwards to apply the indent formatting. Since the bookmark
# paragraph is currently the very last paragraph in the document, its
# (1-based) document index equals the paragraph count before insertion,
# and that is also where the first freshly inserted paragraph will land.
$firstNewIndex = $d.Paragraphs.Count

$insertRange.InsertBefore($blockText)

for ($i = 0; $i -lt $paragraphTexts.Count; $i++) {
    if ($indentedIndexes -contains ($i + 1)) {
        $para = $d.Paragraphs.Item($firstNewIndex + $i)
        $para.Range.ParagraphFormat.LeftIndent = 36
    }
}

# Finally, append the closing "Afterwards" paragraph's text into the
# (still-last) bookmark paragraph, ahead of the bookmark markers themselves.
$bookmarkPara = $d.Paragraphs.Last
$bookmarkPara.Range.InsertBefore("Afterwards: Claire has a concussion and Nicole has a broken collarbone plus a cut just above her left eyebrow. ")

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
